# Update countries & provincias Spain
# Applies the 9-May-2020 data refresh (18:34 -> 19:09) to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Updated timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 19:09"

# --- Updated per-country statistics ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1327919
$ws.Range("C4").Value = 6134
$ws.Range("E4").Value = 1024437
$ws.Range("G4").Value = 234
$ws.Range("H4").Value = 78849

# Row 12: Turquia
$ws.Range("B12").Value = 137115
$ws.Range("C12").Value = 1546
$ws.Range("D12").Value = 89480
$ws.Range("E12").Value = 43896
$ws.Range("F12").Value = 1168
$ws.Range("G12").Value = 50
$ws.Range("H12").Value = 3739

# Row 15: Canada
$ws.Range("B15").Value = 66783
$ws.Range("C15").Value = 349
$ws.Range("D15").Value = 30640
$ws.Range("E15").Value = 31515

# Row 26: Chile
$ws.Range("B26").Value = 27219
$ws.Range("C26").Value = 1247
$ws.Range("D26").Value = 12667
$ws.Range("E26").Value = 14248
$ws.Range("F26").Value = 544
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 304

# Row 28: Irlanda
$ws.Range("B28").Value = 22760
$ws.Range("C28").Value = 219
$ws.Range("E28").Value = 4204
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = 1446

# Row 91: Republica de Yibuti
$ws.Range("B91").Value = 1189
$ws.Range("C91").Value = 54
$ws.Range("D91").Value = 834
$ws.Range("E91").Value = 352

# Row 100: Republica de Chipre
$ws.Range("D100").Value = 401
$ws.Range("E100").Value = 476

# Row 104: Libano
$ws.Range("F104").Value = 4

# Rows 192-193: Belice and Nueva Caledonia swap positions in the ranking
# (Nueva Caledonia's case count overtook Belice's)
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
